# Weekly update: a new price record for "Ciboulette" (Feria Lagunitas de
# Puerto Montt) was collected, so it is inserted as the new top data row
# (row 154, right after the header in row 1 and the existing rows for
# other markets in rows 2-153). All the previously existing rows from 154
# down to 195 shift down by one (to 155-196), keeping their original
# values untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 154; Excel shifts rows 154:195 down to 155:196
# and the new row inherits the formatting (incl. the date number format)
# from the row above it, same as a manual "Insert Row" in the UI.
$ws.Rows.Item(154).Insert()

# Populate the newly inserted row with this week's record.
$ws.Cells.Item(154, 1).Value = 4
$ws.Cells.Item(154, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(154, 3).Value = 'Los Lagos'
$ws.Cells.Item(154, 4).Value = 44642
$ws.Cells.Item(154, 5).Value = 10
$ws.Cells.Item(154, 6).Value = 100112039
$ws.Cells.Item(154, 7).Value = 'Ciboulette'
$ws.Cells.Item(154, 8).Value = 'Sin especificar'
$ws.Cells.Item(154, 9).Value = 'Primera'
$ws.Cells.Item(154, 10).Value = 120
$ws.Cells.Item(154, 11).Value = 6000
$ws.Cells.Item(154, 12).Value = 6000
$ws.Cells.Item(154, 13).Value = 6000
$ws.Cells.Item(154, 14).Value = '$/docena de atados'
$ws.Cells.Item(154, 15).Value = 'Provincia de Cautín'
$ws.Cells.Item(154, 16).Value = 2000
$ws.Cells.Item(154, 17).Value = 3
$ws.Cells.Item(154, 18).Value = 'Hortaliza'
